# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The workbook's column G ("K" = strikeouts, formerly labelled "Strike#")
# holds per-game values that were recomputed by the upstream data pipeline.
# This script rewrites the regenerated K values (column G, rows 2-65) onto
# the existing worksheet, leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values keyed by cell address, as produced by the regenerated
# save_data pipeline (std/mean recalculated, s_vals recomputed).
$kValues = @{
    "G2" = 2
    "G3" = 1
    "G4" = 1
    "G5" = 1
    "G6" = 0
    "G7" = 1
    "G8" = 2
    "G9" = 1
    "G10" = 3
    "G11" = 1
    "G12" = 1
    "G13" = 2
    "G14" = 1
    "G15" = 2
    "G16" = 2
    "G17" = 1
    "G18" = 4
    "G19" = 3
    "G20" = 0
    "G21" = 1
    "G22" = 1
    "G23" = 0
    "G24" = 1
    "G25" = 1
    "G26" = 2
    "G27" = 2
    "G28" = 2
    "G29" = 3
    "G30" = 0
    "G31" = 1
    "G32" = 1
    "G33" = 2
    "G34" = 1
    "G35" = 4
    "G36" = 0
    "G37" = 2
    "G38" = 3
    "G39" = 1
    "G40" = 2
    "G41" = 0
    "G42" = 2
    "G43" = 3
    "G44" = 1
    "G45" = 0
    "G46" = 0
    "G47" = 2
    "G48" = 1
    "G49" = 1
    "G50" = 1
    "G51" = 0
    "G52" = 1
    "G53" = 2
    "G54" = 1
    "G55" = 0
    "G56" = 1
    "G57" = 3
    "G58" = 1
    "G59" = 2
    "G60" = 2
    "G61" = 2
    "G62" = 0
    "G64" = 3
    "G65" = 1
}

foreach ($cellRef in $kValues.Keys) {
    $ws.Range($cellRef).Value = $kValues[$cellRef]
}
